# Update "想去人数" (number of people interested) counts, as published by the
# gh-pages build at commit 456a3b4.
#
# Sheet "展览"   (Exhibitions) - rows 2..22 in column F
# Sheet "演出"   (Performances) - row 2 in column F
# Sheet "全部类型" (All types, aggregate of the above) - rows 2..23 in column F

$wb = $excel.ActiveWorkbook

# Values shared by the "展览" sheet and the first 22 data rows of the
# "全部类型" sheet (row -> new value).
$exhibitionUpdates = @{
    2  = 323
    4  = 10406
    5  = 331
    7  = 19
    8  = 1284
    9  = 7200
    11 = 444
    13 = 128
    14 = 3208
    17 = 672
    18 = 124
    19 = 1038
    21 = 77
    22 = 1645
}

$wsExhibitions = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibitions.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsPerformances = $wb.Worksheets.Item("演出")
$wsPerformances.Range("F2").Value = 29

$wsAll = $wb.Worksheets.Item("全部类型")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsAll.Range("F$row").Value = $exhibitionUpdates[$row]
}
$wsAll.Range("F23").Value = 29
